$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.140252590179443
$ws.Range("B1").Value = 2.690697193145752
$ws.Range("C1").Value = 2.818108081817627
$ws.Range("D1").Value = 2.875888347625732
$ws.Range("E1").Value = 0.7753333449363708
